# ============================================================
# Edit workbook: rebuild "goods_sheet" with new data, reorder
# sheetId, adjust column widths, selections, and active tab.
# ============================================================

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Recreate goods_sheet so it gets a fresh sheetId (5 -> 7), ---
# --- matching the behavior of deleting and re-adding the sheet ---
$oldGoods = $wb.Worksheets.Item("goods_sheet")
$oldGoods.Delete()

$rowdownOrder = $wb.Worksheets.Item("rowdown_order")
$goods = $wb.Worksheets.Add($rowdownOrder)
$goods.Name = "goods_sheet"

$goods.Range("A1").Value = "品名"
$goods.Range("B1").Value = "規格"
$goods.Range("C1").Value = "庫存"
$goods.Range("A2").Value = "蓬萊白米"
$goods.Range("B2").Value = "15kg"
$goods.Range("C2").Value = 2311
$goods.Range("A3").Value = "蓬萊白米"
$goods.Range("B3").Value = "30kg"
$goods.Range("C3").Value = 346
$goods.Range("A4").Value = "長糯白米"
$goods.Range("B4").Value = "30kg"
$goods.Range("C4").Value = 2345
$goods.Range("A5").Value = "圓糯白米"
$goods.Range("B5").Value = "30kg"
$goods.Range("C5").Value = 346
$goods.Range("A6").Value = "清東蓬萊白米"
$goods.Range("B6").Value = "30kg"
$goods.Range("C6").Value = 2354
$goods.Range("A7").Value = "碎米"
$goods.Range("B7").Value = "30kg"
$goods.Range("C7").Value = 4634
$goods.Range("A8").Value = "米糠"
$goods.Range("B8").Value = "30kg"
$goods.Range("C8").Value = 45345
$goods.Range("A9").Value = "長米"
$goods.Range("B9").Value = "30kg"
$goods.Range("C9").Value = 32350
$goods.Range("A10").Value = "蓬萊白米"
$goods.Range("B10").Value = "5kg"
$goods.Range("C10").Value = 111143
$goods.Range("A11").Value = "碎米"
$goods.Range("B11").Value = "3kg"
$goods.Range("C11").Value = 34534
$goods.Range("A12").Value = "長糯白米"
$goods.Range("B12").Value = "10kg"
$goods.Range("C12").Value = 2345
$goods.Range("A13").Value = "清東蓬萊白米"
$goods.Range("B13").Value = "10kg"
$goods.Range("C13").Value = 57547
$goods.Range("A14").Value = "圓糯白米"
$goods.Range("B14").Value = "10kg"
$goods.Range("C14").Value = 346
$goods.Range("A15").Value = "圓糯白米"
$goods.Range("B15").Value = "3kg"
$goods.Range("C15").Value = 6634
$goods.Range("A16").Value = "長米"
$goods.Range("B16").Value = "3kg"
$goods.Range("C16").Value = 32350
$goods.Range("A17").Value = "米糠"
$goods.Range("B17").Value = "3kg"
$goods.Range("C17").Value = 2115
$goods.Range("A18").Value = "清東蓬萊白米"
$goods.Range("B18").Value = "5kg"
$goods.Range("C18").Value = 57547
$goods.Range("A19").Value = "圓糯白米"
$goods.Range("B19").Value = "1kg"
$goods.Range("C19").Value = 346
$goods.Range("A20").Value = "圓糯白米"
$goods.Range("B20").Value = "1.5kg"
$goods.Range("C20").Value = 346
$goods.Range("A21").Value = "圓糯白米"
$goods.Range("B21").Value = "2.5kg"
$goods.Range("C21").Value = 346
$goods.Range("A22").Value = "圓糯白米"
$goods.Range("B22").Value = "5kg"
$goods.Range("C22").Value = 346
$goods.Range("A23").Value = "圓糯白米"
$goods.Range("B23").Value = "6kg"
$goods.Range("C23").Value = 346
$goods.Range("A24").Value = "圓糯白米"
$goods.Range("B24").Value = "7.5kg"
$goods.Range("C24").Value = 346
$goods.Range("A25").Value = "圓糯白米"
$goods.Range("B25").Value = "15kg"
$goods.Range("C25").Value = 346
$goods.Range("A26").Value = "長糯白米"
$goods.Range("B26").Value = "6kg"
$goods.Range("C26").Value = 345
$goods.Range("A27").Value = "蓬萊白米1"
$goods.Range("B27").Value = "15kg"
$goods.Range("C27").Value = 2311
$goods.Range("A28").Value = "蓬萊白米1"
$goods.Range("B28").Value = "30kg"
$goods.Range("C28").Value = 346
$goods.Range("A29").Value = "長糯白米1"
$goods.Range("B29").Value = "30kg"
$goods.Range("C29").Value = 2345
$goods.Range("A30").Value = "圓糯白米1"
$goods.Range("B30").Value = "30kg"
$goods.Range("C30").Value = 346
$goods.Range("A31").Value = "清東蓬萊白米1"
$goods.Range("B31").Value = "30kg"
$goods.Range("C31").Value = 2354
$goods.Range("A32").Value = "碎米1"
$goods.Range("B32").Value = "30kg"
$goods.Range("C32").Value = 4634
$goods.Range("A33").Value = "米糠1"
$goods.Range("B33").Value = "30kg"
$goods.Range("C33").Value = 45345
$goods.Range("A34").Value = "長米1"
$goods.Range("B34").Value = "30kg"
$goods.Range("C34").Value = 32350
$goods.Range("A35").Value = "蓬萊白米1"
$goods.Range("B35").Value = "5kg"
$goods.Range("C35").Value = 111143
$goods.Range("A36").Value = "碎米1"
$goods.Range("B36").Value = "3kg"
$goods.Range("C36").Value = 34534
$goods.Range("A37").Value = "長糯白米1"
$goods.Range("B37").Value = "10kg"
$goods.Range("C37").Value = 2345
$goods.Range("A38").Value = "清東蓬萊白米1"
$goods.Range("B38").Value = "10kg"
$goods.Range("C38").Value = 57547
$goods.Range("A39").Value = "圓糯白米1"
$goods.Range("B39").Value = "10kg"
$goods.Range("C39").Value = 346
$goods.Range("A40").Value = "圓糯白米1"
$goods.Range("B40").Value = "3kg"
$goods.Range("C40").Value = 6634
$goods.Range("A41").Value = "長米1"
$goods.Range("B41").Value = "3kg"
$goods.Range("C41").Value = 32350
$goods.Range("A42").Value = "米糠1"
$goods.Range("B42").Value = "3kg"
$goods.Range("C42").Value = 2115
$goods.Range("A43").Value = "清東蓬萊白米1"
$goods.Range("B43").Value = "5kg"
$goods.Range("C43").Value = 57547
$goods.Range("A44").Value = "長糯白米1"
$goods.Range("B44").Value = "6kg"
$goods.Range("C44").Value = 345
$goods.Range("A45").Value = "蓬萊白米2"
$goods.Range("B45").Value = "15kg"
$goods.Range("C45").Value = 2311
$goods.Range("A46").Value = "蓬萊白米2"
$goods.Range("B46").Value = "30kg"
$goods.Range("C46").Value = 346
$goods.Range("A47").Value = "長糯白米2"
$goods.Range("B47").Value = "30kg"
$goods.Range("C47").Value = 2345
$goods.Range("A48").Value = "圓糯白米2"
$goods.Range("B48").Value = "30kg"
$goods.Range("C48").Value = 346
$goods.Range("A49").Value = "清東蓬萊白米2"
$goods.Range("B49").Value = "30kg"
$goods.Range("C49").Value = 2354
$goods.Range("A50").Value = "碎米2"
$goods.Range("B50").Value = "30kg"
$goods.Range("C50").Value = 4634
$goods.Range("A51").Value = "米糠2"
$goods.Range("B51").Value = "30kg"
$goods.Range("C51").Value = 45345
$goods.Range("A52").Value = "長米2"
$goods.Range("B52").Value = "30kg"
$goods.Range("C52").Value = 32350
$goods.Range("A53").Value = "蓬萊白米2"
$goods.Range("B53").Value = "5kg"
$goods.Range("C53").Value = 111143
$goods.Range("A54").Value = "碎米2"
$goods.Range("B54").Value = "3kg"
$goods.Range("C54").Value = 34534
$goods.Range("A55").Value = "長糯白米2"
$goods.Range("B55").Value = "10kg"
$goods.Range("C55").Value = 2345
$goods.Range("A56").Value = "清東蓬萊白米2"
$goods.Range("B56").Value = "10kg"
$goods.Range("C56").Value = 57547
$goods.Range("A57").Value = "圓糯白米2"
$goods.Range("B57").Value = "10kg"
$goods.Range("C57").Value = 346
$goods.Range("A58").Value = "圓糯白米2"
$goods.Range("B58").Value = "3kg"
$goods.Range("C58").Value = 6634
$goods.Range("A59").Value = "長米2"
$goods.Range("B59").Value = "3kg"
$goods.Range("C59").Value = 32350
$goods.Range("A60").Value = "米糠2"
$goods.Range("B60").Value = "3kg"
$goods.Range("C60").Value = 2115
$goods.Range("A61").Value = "清東蓬萊白米2"
$goods.Range("B61").Value = "5kg"
$goods.Range("C61").Value = 57547
$goods.Range("A62").Value = "長糯白米2"
$goods.Range("B62").Value = "6kg"
$goods.Range("C62").Value = 345

# --- Column widths (closest achievable values via ColumnWidth,
#     which quantizes internally to 1/7-character steps) ---
$goods.Columns.Item(1).ColumnWidth = 22
$goods.Columns.Item(2).ColumnWidth = 17.714285714285715
$goods.Columns.Item(3).ColumnWidth = 19.857142857142858

# --- Selection on goods_sheet ---
$goods.Range("F65").Select()

# --- Update selection on rowdown_order sheet ---
$rowdownOrder2 = $wb.Worksheets.Item("rowdown_order")
$rowdownOrder2.Activate()
$rowdownOrder2.Range("G15").Select()
